# Apply edits described by the commit: rename user, renumber screenshots,
# rewrite explanation text, and move the Windows-update error block from
# rows 5-7 down by shifting it to rows 5 (error) / 7 (operation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: user_name, same replacement for every data row (2..16) ---
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = "Manabu Watanabe"
}

# --- Column J: capimg screenshot path renumbering ---
$capimg = @{
    2  = "bdot20240415_141954/1.png"
    3  = "bdot20240415_141954/2.png"
    4  = "bdot20240415_141954/3.png"
    5  = "bdot20240415_141954/4.png"
    6  = "bdot20240415_141954/5.png"
    7  = "bdot20240415_141954/5.png"
    8  = "bdot20240415_141954/6.png"
    9  = "bdot20240415_141954/7.png"
    10 = "bdot20240415_141954/8.png"
    11 = "bdot20240415_141954/9.png"
    12 = "bdot20240415_141954/10.png"
    13 = "bdot20240415_141954/1.png"
    14 = "bdot20240415_141954/2.png"
    15 = "bdot20240415_141954/3.png"
    16 = "bdot20240415_141954/11.png"
}
foreach ($r in $capimg.Keys) {
    $ws.Cells.Item($r, 10).Value = $capimg[$r]
}

# --- Column K: explanation text ---
$explanation = @{
    2  = "「スタート」ボタンをクリックする"
    3  = "メニューから「設定」アイコンをクリックする"
    4  = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
    5  = "0x80240fff エラー"
    6  = "デスクトップ画面の左下にある「スタート」ボタンを右クリックする"
    7  = "メニューからターミナル(管理者)をクリックする"
    8  = "ユーザーアカウント制御と表示されているウィンドウが開いたことを確認する"
    9  = "PowerShellウィンドウに start-transcript と入力し、[Enter]キーを押す"
    10 = "wuauclt.exe /resetauthorization /detectnow と入力し、[Enter]キーを押す"
    11 = "netsh winhttp show proxy と入力し、[Enter]キーを押す"
    12 = "netsh winhttp reset proxy と入力し、[Enter]キーを押す"
    13 = "「スタート」ボタンをクリックする"
    14 = "メニューから「設定」アイコンをクリックする"
    15 = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
    16 = "「更新プログラムのチェック」ボタンをクリックする"
}
foreach ($r in $explanation.Keys) {
    $ws.Cells.Item($r, 11).Value = $explanation[$r]
}

# --- Row 5: type becomes "error", and error_type / error_content get filled in ---
$ws.Cells.Item(5, 2).Value = "error"
$ws.Cells.Item(5, 12).Value = "Error W"
$ws.Cells.Item(5, 13).Value = " エラーの Windows"

# --- Row 7: type becomes "operation", and error_type / error_content get cleared ---
$ws.Cells.Item(7, 2).Value = "operation"
$ws.Cells.Item(7, 12).Value = ""
$ws.Cells.Item(7, 13).Value = ""
